$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: add A27 (styled marker cell, same fill style as other rows' A column),
# and fill in C27/D27 with new screenshot filenames.
$ws.Range("A26").Copy($ws.Range("A27"))
$ws.Range("C27").Value = "Tests.TestForPositiveInteger.png"
$ws.Range("D27").Value = "Tests that the mean is positive and std dev is above zero.png"

# Update the active selection to match the saved view state.
$ws.Range("B33").Select()
